# Time Analysis.xlsx edit: "Data for 1000 column"
# 1) Remove the old "Randomness Sort" block (rows 34-38, the placeholder
#    group with stray ~8000 values) by deleting those entire rows. This
#    shifts every following row group up by 5 and drops the now-unused
#    "Randomness Sort" shared string automatically on save.
# 2) Fill in the newly-added "1000" input-size column (E) with trial/avg
#    timing data for every sort's row group (rows 4-63 after the shift).
# 3) Restore the sheet's selection/active cell to match the saved state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A34:A38").EntireRow.Delete()

$ws.Range("E4").Value2 = 14
$ws.Range("E5").Value2 = 12
$ws.Range("E6").Value2 = 5
$ws.Range("E7").Value2 = 4
$ws.Range("E8").Value2 = 5
$ws.Range("E9").Value2 = 17
$ws.Range("E10").Value2 = 13
$ws.Range("E11").Value2 = 4
$ws.Range("E12").Value2 = 3
$ws.Range("E13").Value2 = 2
$ws.Range("E14").Value2 = 21
$ws.Range("E15").Value2 = 6
$ws.Range("E16").Value2 = 4
$ws.Range("E17").Value2 = 2
$ws.Range("E18").Value2 = 2
$ws.Range("E19").Value2 = 4
$ws.Range("E20").Value2 = 1
$ws.Range("E21").Value2 = 1
$ws.Range("E22").Value2 = 1
$ws.Range("E23").Value2 = 1
$ws.Range("E24").Value2 = 33
$ws.Range("E25").Value2 = 69
$ws.Range("E26").Value2 = 62
$ws.Range("E27").Value2 = 2
$ws.Range("E28").Value2 = 1
$ws.Range("E29").Value2 = 2
$ws.Range("E30").Value2 = 3
$ws.Range("E31").Value2 = 1
$ws.Range("E32").Value2 = 2
$ws.Range("E33").Value2 = 2
$ws.Range("E34").Value2 = 12
$ws.Range("E35").Value2 = 6
$ws.Range("E36").Value2 = 9
$ws.Range("E37").Value2 = 3
$ws.Range("E38").Value2 = 4
$ws.Range("E39").Value2 = 1
$ws.Range("E40").Value2 = 1
$ws.Range("E41").Value2 = 1
$ws.Range("E42").Value2 = 1
$ws.Range("E43").Value2 = 1
$ws.Range("E44").Value2 = 1
$ws.Range("E45").Value2 = 1
$ws.Range("E46").Value2 = 1
$ws.Range("E47").Value2 = 1
$ws.Range("E48").Value2 = 1
$ws.Range("E49").Value2 = 2
$ws.Range("E50").Value2 = 2
$ws.Range("E51").Value2 = 3
$ws.Range("E52").Value2 = 2
$ws.Range("E53").Value2 = 1
$ws.Range("E54").Value2 = 2
$ws.Range("E55").Value2 = 1
$ws.Range("E56").Value2 = 1
$ws.Range("E57").Value2 = 1
$ws.Range("E58").Value2 = 0
$ws.Range("E59").Value2 = 1
$ws.Range("E60").Value2 = 1
$ws.Range("E61").Value2 = 1
$ws.Range("E62").Value2 = 0
$ws.Range("E63").Value2 = 1

$ws.Range("E64").Select()
